$wb = $excel.ActiveWorkbook

# --- "partnership" sheet: rename headers and update the share series ---
$ws = $wb.Worksheets.Item("partnership")

$ws.Range("A1").Value = "Year"
$ws.Range("B1").Value = "Share"

$values = @(
    0.60599999999999998,
    0.60860000000000003,
    0.60570000000000002,
    0.61019999999999996,
    0.61339999999999995,
    0.61570000000000003,
    0.61919999999999997,
    0.62190000000000001,
    0.624,
    0.62809999999999999,
    0.64049999999999996,
    0.62729999999999997,
    0.62050000000000005,
    0.61850000000000005
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
